# Commit: created skeleton lakeidGenerator and added S3 key integration
# (S3 not integrated yet, just defined keys) -- spreadsheet-side change:
# add a "Column Comparison" sheet mapping CSV headers to Excel headers,
# inserted between the CSV-data sheet and the Description sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("GLEON_GMA_Example csv")
$wsDescription = $wb.Worksheets.Item("Description")

# --- create the new sheet and place it between the two existing sheets ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Column Comparison"
$ws1 = $wb.Worksheets.Item("GLEON_GMA_Example csv")
$newSheet.Move($null, $ws1)

$cc = $wb.Worksheets.Item("Column Comparison")

# --- column widths ---
$cc.Columns.Item(1).ColumnWidth = 23.1640625
$cc.Columns.Item(2).ColumnWidth = 23.1640625

# --- data: CSV Columns / Excel Columns comparison table ---
$cc.Range("A1").Value = 'CSV Columns'
$cc.Range("B1").Value = 'Excel Columns'
$cc.Range("A2").Value = 'Date'
$cc.Range("B2").Value = 'Lake_ID (generated by us)'
$cc.Range("A3").Value = 'Lat'
$cc.Range("B3").Value = 'Date'
$cc.Range("A4").Value = 'Long'
$cc.Range("B4").Value = 'LakeName'
$cc.Range("A5").Value = 'Microcystin'
$cc.Range("B5").Value = 'DataContact'
$cc.Range("A6").Value = 'Anatoxin-a'
$cc.Range("B6").Value = 'Lat'
$cc.Range("A7").Value = 'Cylindrospermopsin'
$cc.Range("B7").Value = 'Long'
$cc.Range("A8").Value = 'Saxitoxin'
$cc.Range("B8").Value = 'Altitude_m'
$cc.Range("A9").Value = 'Nodularin'
$cc.Range("B9").Value = 'MaximumDepth_m'
$cc.Range("A10").Value = 'Geosmin'
$cc.Range("B10").Value = 'MeanDepth_m'
$cc.Range("A11").Value = '2-MIB'
$cc.Range("B11").Value = 'SecchiDepth_m'
$cc.Range("A12").Value = 'TN'
$cc.Range("B12").Value = 'SamplingDepth_m'
$cc.Range("A13").Value = 'TP'
$cc.Range("B13").Value = 'ThermoclineDepth_m'
$cc.Range("A14").Value = 'Secchi'
$cc.Range("B14").Value = 'SurfaceTemperature_C'
$cc.Range("A15").Value = 'Chl'
$cc.Range("B15").Value = 'EpilimneticTemperature_C'
$cc.Range("A16").Value = 'NO2+3'
$cc.Range("B16").Value = 'TN_mgL'
$cc.Range("A17").Value = 'NO3'
$cc.Range("B17").Value = 'TP_mgL'
$cc.Range("A18").Value = 'NH3'
$cc.Range("B18").Value = 'NO3NO2_mgL'
$cc.Range("A19").Value = 'OrthoP'
$cc.Range("B19").Value = 'NH4_mgL'
$cc.Range("A20").Value = 'SRP'
$cc.Range("B20").Value = 'PO4_ugL'
$cc.Range("A21").Value = 'TotalPhytoCells'
$cc.Range("B21").Value = 'Chlorophylla_ugL'
$cc.Range("A22").Value = 'CyanobacterialCells'
$cc.Range("B22").Value = 'Chlorophyllb_ugL'
$cc.Range("A23").Value = 'PercentCyano'
$cc.Range("B23").Value = 'Zeaxanthin_ugL'
$cc.Range("A24").Value = 'DominantBloomGenera'
$cc.Range("B24").Value = 'Diadinoxanthin_ugL'
$cc.Range("A25").Value = 'mcyDgeneAbund'
$cc.Range("B25").Value = 'Fucoxanthin_ugL'
$cc.Range("A26").Value = 'mcyEgeneAbund'
$cc.Range("B26").Value = 'Diatoxanthin_ugL'
$cc.Range("A27").Value = 'Comments'
$cc.Range("B27").Value = 'Alloxanthin_ugL'
$cc.Range("B28").Value = 'Peridinin_ugL'
$cc.Range("B29").Value = 'Chlorophyllc2_ugL'
$cc.Range("B30").Value = 'Echinenone_ugL'
$cc.Range("B31").Value = 'Lutein_ugL'
$cc.Range("B32").Value = 'Violaxanthin_ugL'
$cc.Range("B33").Value = 'TotalMC_ug/L'
$cc.Range("B34").Value = 'DissolvedMC_ugL'
$cc.Range("B35").Value = 'MC_YR_ugL'
$cc.Range("B36").Value = 'MC_dmRR_ugL'
$cc.Range("B37").Value = 'MC_RR_ugL'
$cc.Range("B38").Value = 'MC_dmLR_ugL'
$cc.Range("B39").Value = 'MC_LR_ugL'
$cc.Range("B40").Value = 'MC_LY_ugL'
$cc.Range("B41").Value = 'MC_LW_ugL'
$cc.Range("B42").Value = 'MC_LF_ugL'
$cc.Range("B43").Value = 'NOD_ugL'
$cc.Range("B44").Value = 'CYN_ugL'
$cc.Range("B45").Value = 'ATX_ugL'
$cc.Range("B46").Value = 'SAX_ugL'
$cc.Range("B47").Value = 'GEO_ngL'
$cc.Range("B48").Value = '2MIB_ngL'
$cc.Range("B49").Value = 'TotalPhyto_CellsmL'
$cc.Range("B50").Value = 'Cyano_CellsmL'
$cc.Range("B51").Value = 'PercentCyano'
$cc.Range("B52").Value = 'DominantBloomGenera'
$cc.Range("B53").Value = 'mcyD_genemL'
$cc.Range("B54").Value = 'mcyE_genemL'
$cc.Range("B55").Value = 'Comments'

# --- formatting: reuse the existing bold header styles already present
#     on the "GLEON_GMA_Example csv" sheet (style "2" = Calibri bold used
#     on A1:F1, style "1" = Source Sans Pro bold used on G1:BB1) so no new
#     fonts/styles get introduced into styles.xml ---
$ws1.Range("A1").Copy()
$cc.Range("B2:B7").PasteSpecial(-4122)
$ws1.Range("G1").Copy()
$cc.Range("B8:B32").PasteSpecial(-4122)
$ws1.Range("A1").Copy()
$cc.Range("B33:B34").PasteSpecial(-4122)
$ws1.Range("G1").Copy()
$cc.Range("B35:B45").PasteSpecial(-4122)
$ws1.Range("A1").Copy()
$cc.Range("B46:B55").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- selection / active sheet bookkeeping ---
$ws1.Select()
$ws1.Range("A4").Select()

$cc.Select()
$cc.Range("F17").Select()

$wsDescription.Select()
$wsDescription.Range("B8").Select()

Write-Host "done"
